$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: typo fix SDS -> SSD
$ws.Range("D5").Value = "SSD"

# Gray-out old rows 11, 21, 22 with the new "light gray" style
$ws.Range("A11:D11").Font.ThemeColor = 2
$ws.Range("A11:D11").Font.TintAndShade = -0.249977111117893
$ws.Range("A21:D21").Font.ThemeColor = 2
$ws.Range("A21:D21").Font.TintAndShade = -0.249977111117893
$ws.Range("A22:D22").Font.ThemeColor = 2
$ws.Range("A22:D22").Font.TintAndShade = -0.249977111117893

# New rows 29-30
$ws.Range("B29").Value = "exact replications"
$ws.Range("C29").Value = "only linear regression"
$ws.Range("D29").Value = "exact"

$ws.Range("B30").Value = "conceptual replications"
$ws.Range("C30").Value = "linear, logistic, probit regression"
$ws.Range("D30").Value = "concept"

# Sheet view changes
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F36:G36").Select()
